# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh values to the Leve profit sheets
# (columns H-N: currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 3921.9375
$ws.Range("J103").Value = 767.6429000000001
$ws.Range("L103").Value = 2302.9287
$ws.Range("N103").Value = -3474.9287
$ws.Range("H120").Value = 30000
$ws.Range("J120").Value = 30000
$ws.Range("L120").Value = 30000
$ws.Range("N120").Value = -39676
$ws.Range("H129").Value = 655.26666
$ws.Range("I129").Value = 414.23077
$ws.Range("J129").Value = 2222
$ws.Range("K129").Value = 1242.69231
$ws.Range("L129").Value = 6666
$ws.Range("M129").Value = 3757.30769
$ws.Range("N129").Value = -16666

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2236.6316
$ws.Range("I61").Value = 2154.8572
$ws.Range("J61").Value = 2465.6
$ws.Range("K61").Value = 2154.8572
$ws.Range("L61").Value = 2465.6
$ws.Range("M61").Value = -1942.8572
$ws.Range("N61").Value = -2889.6
$ws.Range("H74").Value = 39484.168
$ws.Range("I74").Value = 64661.75
$ws.Range("J74").Value = 10709.786
$ws.Range("K74").Value = 64661.75
$ws.Range("L74").Value = 10709.786
$ws.Range("M74").Value = -63787.75
$ws.Range("N74").Value = -12457.786
$ws.Range("H77").Value = 39484.168
$ws.Range("I77").Value = 64661.75
$ws.Range("J77").Value = 10709.786
$ws.Range("K77").Value = 323308.75
$ws.Range("L77").Value = 53548.93
$ws.Range("M77").Value = -318940.75
$ws.Range("N77").Value = -62284.93
$ws.Range("H102").Value = 2355.8333
$ws.Range("I102").Value = 2217.6667
$ws.Range("J102").Value = 2770.3333
$ws.Range("K102").Value = 2217.6667
$ws.Range("L102").Value = 2770.3333
$ws.Range("M102").Value = -595.6667000000002
$ws.Range("N102").Value = -6014.3333
$ws.Range("H136").Value = 2236.6316
$ws.Range("I136").Value = 2154.8572
$ws.Range("J136").Value = 2465.6
$ws.Range("K136").Value = 6464.571599999999
$ws.Range("L136").Value = 7396.799999999999
$ws.Range("M136").Value = -3914.571599999999
$ws.Range("N136").Value = -12496.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 398.75
$ws.Range("I5").Value = 318.18182
$ws.Range("J5").Value = 576
$ws.Range("K5").Value = 318.18182
$ws.Range("L5").Value = 576
$ws.Range("M5").Value = -205.18182
$ws.Range("N5").Value = -802
$ws.Range("H134").Value = 681477.75
$ws.Range("I134").Value = 912439.3
$ws.Range("J134").Value = 3990.4666
$ws.Range("K134").Value = 2737317.9
$ws.Range("L134").Value = 11971.3998
$ws.Range("M134").Value = -2734782.9
$ws.Range("N134").Value = -17041.3998

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 3340
$ws.Range("J4").Value = 3340
$ws.Range("L4").Value = 3340
$ws.Range("N4").Value = -3564
$ws.Range("H31").Value = 1553.34
$ws.Range("I31").Value = 929
$ws.Range("J31").Value = 2286.261
$ws.Range("K31").Value = 929
$ws.Range("L31").Value = 2286.261
$ws.Range("M31").Value = -634
$ws.Range("N31").Value = -2876.261
$ws.Range("H34").Value = 1553.34
$ws.Range("I34").Value = 929
$ws.Range("J34").Value = 2286.261
$ws.Range("K34").Value = 929
$ws.Range("L34").Value = 2286.261
$ws.Range("M34").Value = -727
$ws.Range("N34").Value = -2690.261
$ws.Range("H58").Value = 14256.25
$ws.Range("I58").Value = 21450
$ws.Range("K58").Value = 21450
$ws.Range("M58").Value = -21247
$ws.Range("H99").Value = 3207164
$ws.Range("I99").Value = 4168488
$ws.Range("J99").Value = 2750.6667
$ws.Range("K99").Value = 4168488
$ws.Range("L99").Value = 2750.6667
$ws.Range("M99").Value = -4166990
$ws.Range("N99").Value = -5746.6667
$ws.Range("H107").Value = 2468.182
$ws.Range("I107").Value = 858.3333
$ws.Range("J107").Value = 4400
$ws.Range("K107").Value = 858.3333
$ws.Range("L107").Value = 4400
$ws.Range("M107").Value = 1061.6667
$ws.Range("N107").Value = -8240
$ws.Range("H122").Value = 1416.6666
$ws.Range("I122").Value = 1390
$ws.Range("J122").Value = 1450
$ws.Range("K122").Value = 4170
$ws.Range("L122").Value = 4350
$ws.Range("M122").Value = -1720
$ws.Range("N122").Value = -9250
$ws.Range("H126").Value = 3207164
$ws.Range("I126").Value = 4168488
$ws.Range("J126").Value = 2750.6667
$ws.Range("K126").Value = 12505464
$ws.Range("L126").Value = 8252.000100000001
$ws.Range("M126").Value = -12502994
$ws.Range("N126").Value = -13192.0001
$ws.Range("H132").Value = 743100.0600000001
$ws.Range("I132").Value = 2307.5
$ws.Range("J132").Value = 2317284.2
$ws.Range("K132").Value = 6922.5
$ws.Range("L132").Value = 6951852.600000001
$ws.Range("M132").Value = -4392.5
$ws.Range("N132").Value = -6956912.600000001
$ws.Range("H133").Value = 28105.334
$ws.Range("J133").Value = 28105.334
$ws.Range("L133").Value = 28105.334
$ws.Range("N133").Value = -33165.334
$ws.Range("H136").Value = 14256.25
$ws.Range("I136").Value = 21450
$ws.Range("K136").Value = 64350
$ws.Range("M136").Value = -61800
$ws.Range("H141").Value = 27545.2
$ws.Range("J141").Value = 27545.2
$ws.Range("L141").Value = 27545.2
$ws.Range("N141").Value = -37905.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 1030
$ws.Range("I138").Value = 1030
$ws.Range("K138").Value = 3090
$ws.Range("M138").Value = 2050

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 4248.5
$ws.Range("I4").Value = 2000
$ws.Range("J4").Value = 4998
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 4998
$ws.Range("M4").Value = -1888
$ws.Range("N4").Value = -5222
$ws.Range("H126").Value = 71435096
$ws.Range("I126").Value = 90915850
$ws.Range("J126").Value = 5666.6665
$ws.Range("K126").Value = 272747550
$ws.Range("L126").Value = 16999.9995
$ws.Range("M126").Value = -272745080
$ws.Range("N126").Value = -21939.9995

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 66890
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 66890
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 66890
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -67116
$ws.Range("H7").Value = 4329.533
$ws.Range("I7").Value = 4283.6665
$ws.Range("J7").Value = 4398.3335
$ws.Range("K7").Value = 4283.6665
$ws.Range("L7").Value = 4398.3335
$ws.Range("M7").Value = -4171.6665
$ws.Range("N7").Value = -4622.3335
$ws.Range("H28").Value = 66890
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 66890
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 66890
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -67354
$ws.Range("H37").Value = 66890
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 66890
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 66890
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -67104
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H126").Value = 4329.533
$ws.Range("I126").Value = 4283.6665
$ws.Range("J126").Value = 4398.3335
$ws.Range("K126").Value = 12850.9995
$ws.Range("L126").Value = 13195.0005
$ws.Range("M126").Value = -10380.9995
$ws.Range("N126").Value = -18135.0005

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1799.8
$ws.Range("I14").Value = 999
$ws.Range("K14").Value = 999
$ws.Range("M14").Value = -831
$ws.Range("H126").Value = 1413.5333
$ws.Range("I126").Value = 1306
$ws.Range("J126").Value = 1574.8334
$ws.Range("K126").Value = 3918
$ws.Range("L126").Value = 4724.5002
$ws.Range("M126").Value = -1448
$ws.Range("N126").Value = -9664.5002
